$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit rotates the species-record data (columns A,B,D,E,F,G,H,P,Q,R)
# among rows 48, 50, 51, 52 and 53, while leaving all other columns
# (C, I, J, S, T, U, V, W, Y, Z, AA, AB, AD, AE, AG, AI, AT, AW, AX, AY)
# untouched for each physical row.
#
# Data flow (new row <- old row):
#   48 <- 52
#   50 <- 53
#   51 <- 50
#   52 <- 51
#   53 <- 48

$cols = @("A","B","D","E","F","G","H","P","Q","R")

function Get-RowData($row) {
    $data = @{}
    foreach ($col in $cols) {
        # Note: reading the ".Value" getter through this COM-interop layer
        # returns the property descriptor rather than the actual value, so
        # ".Value2" must be used instead for reading cell contents.
        $data[$col] = $ws.Range("$col$row").Value2
    }
    return $data
}

# Snapshot the original values of the affected rows before overwriting anything.
$row48 = Get-RowData 48
$row50 = Get-RowData 50
$row51 = Get-RowData 51
$row52 = Get-RowData 52
$row53 = Get-RowData 53

function Set-RowData($row, $data) {
    foreach ($col in $cols) {
        $ws.Range("$col$row").Value2 = $data[$col]
    }
}

Set-RowData 48 $row52
Set-RowData 50 $row53
Set-RowData 51 $row50
Set-RowData 52 $row51
Set-RowData 53 $row48
